# Insert a new weekly record at the top of the price history (row 14) and
# push every existing record (rows 14-65) down by one row, creating a new
# row 66 for what used to be the last record (old row 65).
#
# Only columns D (Fecha), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg) vary
# row to row - all the other columns (A,B,C,E,F,G,H,N,O,Q,R) are constant
# across the whole sheet, so they're simply copied verbatim for the
# brand-new row 66.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make the new row 66 exist with the same "constant" columns as every other
# row (copy the full row 65 as a template), and the same date style/format
# on column D.
$template = $ws.Range("A65:R65").Value2
$ws.Range("A66:R66").Value2 = $template
$ws.Range("D66").NumberFormat = $ws.Range("D65").NumberFormat

# Shift the variable columns down by one row, from the bottom up so we
# never overwrite a source row before it has been read.
for ($row = 66; $row -ge 15; $row--) {
    $srcRow = $row - 1
    $ws.Cells.Item($row, 4).Value2  = $ws.Cells.Item($srcRow, 4).Value2   # D Fecha
    $ws.Cells.Item($row, 9).Value2  = $ws.Cells.Item($srcRow, 9).Value2   # I Calidad
    $ws.Cells.Item($row, 10).Value2 = $ws.Cells.Item($srcRow, 10).Value2  # J Volumen
    $ws.Cells.Item($row, 11).Value2 = $ws.Cells.Item($srcRow, 11).Value2  # K Precio minimo
    $ws.Cells.Item($row, 12).Value2 = $ws.Cells.Item($srcRow, 12).Value2  # L Precio maximo
    $ws.Cells.Item($row, 13).Value2 = $ws.Cells.Item($srcRow, 13).Value2  # M Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value2 = $ws.Cells.Item($srcRow, 16).Value2  # P Precio $/Kg
}

# Row 14 now holds the new weekly record. Calidad (I14) and Volumen (J14)
# keep their previous values ("Primera" / 100); only the date and the three
# price columns plus $/Kg change.
$ws.Cells.Item(14, 4).Value2  = 45037
$ws.Cells.Item(14, 11).Value2 = 12500
$ws.Cells.Item(14, 12).Value2 = 13000
$ws.Cells.Item(14, 13).Value2 = 12750
$ws.Cells.Item(14, 16).Value2 = 708
